$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 90) with the 2023-06-29 11:43:31 resale-number
# update, matching the shape of all the existing rows.
#
# Columns A (date) and D (week number) hold values that *look* numeric/date
# to Excel's auto-detection, which would otherwise silently convert them to
# a date serial / a plain number (and pick up a stray number-format style in
# the process). To keep them as genuine text - exactly like every other row
# in the sheet - we build them via a text formula and then "paste values"
# over themselves, which locks in the text result without reformatting the
# cell or leaving orphaned styles behind.
$ws.Range("A90").Formula = '="2023-06-29"'
$ws.Range("A90").Copy()
$ws.Range("A90").PasteSpecial(-4163)

$ws.Range("B90").Value = "11:43:31"
$ws.Range("C90").Value = "Thursday"

$ws.Range("D90").Formula = '="26"'
$ws.Range("D90").Copy()
$ws.Range("D90").PasteSpecial(-4163)

$ws.Range("E90").Value = 123159
$ws.Range("F90").Value = 134486
$ws.Range("G90").Value = 163790
$ws.Range("H90").Value = 134250
$ws.Range("I90").Value = 177130
$ws.Range("J90").Value = 115115
$ws.Range("K90").Value = 204447
$ws.Range("L90").Value = 226420
$ws.Range("M90").Value = 176344
$ws.Range("N90").Value = 104510
$ws.Range("O90").Value = 39766
$ws.Range("P90").Value = 33724
$ws.Range("Q90").Value = 52519
$ws.Range("R90").Value = -1
$ws.Range("S90").Value = 35760
$ws.Range("T90").Value = -1
